{"js": "// Locate the \"Size axis captions\" list item and:\n//   1) change its text to \"Size axis - captions\"\n//   2) add a new list item right after it: \"Set x,y-ticks\"\nconst body = context.document.body;\nconst results = body.search(\"Size axis captions\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"Size axis captions\" in the document.');\n}\n\nconst target = results.items[0];\n\n// Replace the text in place so the paragraph keeps its original\n// formatting (style + Times New Roman run font).\ntarget.insertText(\"Size axis - captions\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Insert a brand-new list paragraph right after it, inheriting the\n// same \"List Paragraph\" style / font used throughout this list.\nconst para = target.paragraphs.getFirst();\npara.insertParagraph(\"Set x,y-ticks\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Locate the \"Size axis captions\" list item and:\n#   1) change its text to \"Size axis - captions\"\n#   2) add a new list item right after it: \"Set x,y-ticks\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Size axis captions\"\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif (-not $found) {\n    throw 'Could not find \"Size axis captions\" in the document.'\n}\n\n$rng = $find.Parent\n$rng.Text = \"Size axis - captions\"\n\n# Insert a brand-new list paragraph right after it, inheriting the same\n# \"List Paragraph\" style / Times New Roman font used throughout this list.\n$para = $rng.Paragraphs(1)\n$para.Range.InsertParagraphAfter()\n$nextPara = $para.Next()\n$nextPara.Range.Text = \"Set x,y-ticks\"\n"}
